$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 0.2100747399255385
$ws.Range("D2").Value = 5888186138.283437
$ws.Range("G2").Value = 0.3
$ws.Range("I2").Value = 2850696028.994353
$ws.Range("M2").Value = 11159709000

# Row 3 updates
$ws.Range("B3").Value = 0.2100747399255385
$ws.Range("D3").Value = 5888186138.283437
$ws.Range("G3").Value = 0.3
$ws.Range("I3").Value = 2850696028.994353
$ws.Range("M3").Value = 11159709000
